$p = $ppt.ActivePresentation
$dt = $p.DocumentTheme
$m = $dt | Get-Member
Write-Output $m
